# Final Logistic Regression Model trained
# The "Eye"/"Shoe" category sample rows (rows 4493:4618, referencing the
# now-removed "Vision Express" ... "Guru kirpa custom Shop" shared strings)
# are cleared out of the Name/Category tagging sheet - the cells are wiped
# but the (now blank) rows are left in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4493:B4618").ClearContents()

# Leave the selection where the author ended up after deleting the block.
$ws.Range("C4488").Select()
